# Add 2022-Q3 data
# 1) Insert a new row at the top of the "总计" (summary) sheet's data block
#    with the 2022-Q3 totals.
# 2) Insert a brand-new worksheet named "2022-Q3" right after "总计" (i.e.
#    before "2022-Q2"), populated with the per-fund holdings table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - insert new row 2 with the 2022-Q3 summary values
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Copy formatting (bold/border/center) from the row below onto the new A2
# so it matches the other index cells in column A.
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)

$summary.Cells.Item(2,1).Value2 = 0
$summary.Cells.Item(2,2).Value2 = "2022-Q3"
$summary.Cells.Item(2,3).Value2 = 4
$summary.Cells.Item(2,4).Value2 = 0.43

# The freshly-inserted row inherits the header's bold style for B:D - reset
# those back to the plain/default style used by the rest of the data rows.
$summary.Range("B2:D2").Style = "Normal"

# Column A is a running 0-based index - renumber every data row now that a
# new one was inserted at the top (row 3 was "0", now it must read "1", ...).
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value2 = $r - 2
}

# ---------------------------------------------------------------------
# Step 2: add a new "2022-Q3" worksheet positioned right after "总计"
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q3"
$newSheet.Move($anchor)

# Re-fetch by name: the old $newSheet reference is stale once moved.
$q3 = $wb.Worksheets.Item("2022-Q3")

# Make the new sheet's outline properties match the rest of the workbook's
# sheets (<sheetPr><outlinePr .../><pageSetUpPr/></sheetPr>).
$q3.Outline.SummaryBelow = 1
$q3.Outline.SummaryRight = 1

# -- Header row (bold / centered / bordered, same visual style as s=2) --
$summary.Range("B1:D1").Copy()
$q3.Range("B1").PasteSpecial(-4122)
$q3.Cells.Item(1,2).Value2 = "基金代码"
$q3.Cells.Item(1,3).Value2 = "基金名称"
$q3.Cells.Item(1,4).Value2 = "基金规模"
$q3.Cells.Item(1,5).Value2 = "股票总仓位"
$q3.Cells.Item(1,6).Value2 = "仓位占比"
$q3.Cells.Item(1,7).Value2 = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value2 = "仓位排名"
$q3.Range("E1:H1").Style = $q3.Range("B1").Style

# -- Data rows --
# Columns B..G are text (fund codes / names / numeric-looking strings kept
# as text), column A (row index) and H (rank) are real numbers.
$rows = @(
    @(0, "011466", "兴业医疗保健混合A", "3.85", "87.15", "6.46", "0.2487", 2),
    @(1, "011467", "兴业医疗保健混合C", "1.94", "87.15", "6.46", "0.1253", 2),
    @(2, "008619", "永赢医药健康股票C", "0.40", "94.40", "8.10", "0.0324", 4),
    @(3, "008618", "永赢医药健康股票A", "0.24", "94.40", "8.10", "0.0194", 4)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q3.Cells.Item($r, 1).Value2 = $data[0]

    $textRange = $q3.Range("B$r`:G$r")
    $textRange.NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value2 = $data[1]
    $q3.Cells.Item($r, 3).Value2 = $data[2]
    $q3.Cells.Item($r, 4).Value2 = $data[3]
    $q3.Cells.Item($r, 5).Value2 = $data[4]
    $q3.Cells.Item($r, 6).Value2 = $data[5]
    $q3.Cells.Item($r, 7).Value2 = $data[6]
    $textRange.Style = "Normal"

    $q3.Cells.Item($r, 8).Value2 = $data[7]
}
